$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.603237271308899
$ws.Range("B1").Value = 1.58710241317749
$ws.Range("C1").Value = 1.567187786102295
$ws.Range("D1").Value = 1.956164598464966
$ws.Range("E1").Value = 2.961583614349365
